$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "82 x 17" + [char]11 + "  1    7" + [char]11 + "  ----" + [char]11 + "8|    |" + [char]11 + "2|    |"
$t.Cell(1,2).Range.Text = "23 x 81" + [char]11 + "  8    1" + [char]11 + "  ----" + [char]11 + "2|    |" + [char]11 + "3|    |"
$t.Cell(1,3).Range.Text = "88 x 50" + [char]11 + "  5    0" + [char]11 + "  ----" + [char]11 + "8|    |" + [char]11 + "8|    |"
$t.Cell(2,1).Range.Text = "97 x 45" + [char]11 + "  4    5" + [char]11 + "  ----" + [char]11 + "9|    |" + [char]11 + "7|    |"
$t.Cell(2,2).Range.Text = "98 x 64" + [char]11 + "  6    4" + [char]11 + "  ----" + [char]11 + "9|    |" + [char]11 + "8|    |"
$t.Cell(2,3).Range.Text = "19 x 28" + [char]11 + "  2    8" + [char]11 + "  ----" + [char]11 + "1|    |" + [char]11 + "9|    |"
$t.Cell(3,1).Range.Text = "66 x 55" + [char]11 + "  5    5" + [char]11 + "  ----" + [char]11 + "6|    |" + [char]11 + "6|    |"
$t.Cell(3,2).Range.Text = "99 x 66" + [char]11 + "  6    6" + [char]11 + "  ----" + [char]11 + "9|    |" + [char]11 + "9|    |"
$t.Cell(3,3).Range.Text = "79 x 30" + [char]11 + "  3    0" + [char]11 + "  ----" + [char]11 + "7|    |" + [char]11 + "9|    |"
$t.Cell(4,1).Range.Text = "70 x 14" + [char]11 + "  1    4" + [char]11 + "  ----" + [char]11 + "7|    |" + [char]11 + "0|    |"
$t.Cell(4,2).Range.Text = "40 x 59" + [char]11 + "  5    9" + [char]11 + "  ----" + [char]11 + "4|    |" + [char]11 + "0|    |"
$t.Cell(4,3).Range.Text = "43 x 49" + [char]11 + "  4    9" + [char]11 + "  ----" + [char]11 + "4|    |" + [char]11 + "3|    |"
$t.Cell(5,1).Range.Text = "48 x 30" + [char]11 + "  3    0" + [char]11 + "  ----" + [char]11 + "4|    |" + [char]11 + "8|    |"
$t.Cell(5,2).Range.Text = "43 x 37" + [char]11 + "  3    7" + [char]11 + "  ----" + [char]11 + "4|    |" + [char]11 + "3|    |"
$t.Cell(5,3).Range.Text = "84 x 78" + [char]11 + "  7    8" + [char]11 + "  ----" + [char]11 + "8|    |" + [char]11 + "4|    |"
